$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.803.23"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.598.65"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.475"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.245"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").Value = "1.820.96"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "1.595.79"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").Value = "25.809.67"
$ws.Range("E16").Value = "  -0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("E18").Value = "  -3.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("D36").Value = "1.095.97"
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("E37").Value = "  -2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.795"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0151"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.491"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").Value = "1.733.47"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.739"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").Value = "0.0₇0995"
$ws.Range("E46").Value = "  -12.74%  "
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.65%  "
